$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per the source data refresh.
# Price-column values that read as plain numbers must be forced to stay text
# (matching the original inlineStr cells) via a temporary text NumberFormat,
# then the style is reset to Normal so no stray formatting is introduced.

$ws.Range('D2').Value = '43.800.05'
$ws.Range('E2').Value = '  -0.32%  '
$ws.Range('D3').Value = '2.344.78'
$ws.Range('E3').Value = '  -0.96%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.71'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.669'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.88%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '72.35'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -4.81%  '
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.591'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.99%  '
$ws.Range('E10').Value = '  -2.76%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '58.24'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.38%  '
$ws.Range('E12').Value = '  -0.68%  '
$ws.Range('E13').Value = '  -0.26%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.16'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.38%  '
$ws.Range('D15').Value = '2.692.15'
$ws.Range('E15').Value = '  -1.13%  '
$ws.Range('E16').Value = '  -3.73%  '
$ws.Range('E17').Value = '  -2.24%  '
$ws.Range('D18').Value = '2.344.77'
$ws.Range('E18').Value = '  -0.55%  '
$ws.Range('D19').Value = '43.697.97'
$ws.Range('E19').Value = '  -0.64%  '
$ws.Range('E20').Value = '  -0.87%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.66'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '78.39'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.64%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '253.69'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.19%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.92'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +8.35%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.999'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.07%  '
$ws.Range('E26').Value = '  +2.79%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.50'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.20%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.39'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -6.27%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.26'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.89%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '175.80'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.15%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.27'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.92%  '
$ws.Range('E32').Value = '  -1.30%  '
$ws.Range('E33').Value = '  +0.68%  '
$ws.Range('E34').Value = '  -2.16%  '
$ws.Range('E35').Value = '  -4.30%  '
$ws.Range('E36').Value = '  +0.47%  '
$ws.Range('E37').Value = '  -2.22%  '
$ws.Range('E38').Value = '  -3.60%  '
$ws.Range('E39').Value = '  -2.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0271'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.71%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.24'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +16.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '64.67'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +18.34%  '
$ws.Range('E43').Value = '  +2.66%  '
$ws.Range('E44').Value = '  +6.67%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.78'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.02%  '
$ws.Range('E46').Value = '  -1.47%  '
$ws.Range('E47').Value = '  -0.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.45'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.53%  '
$ws.Range('E49').Value = '  -3.09%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '98.28'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.70%  '
$ws.Range('E51').Value = '  -4.68%  '
